# Re-run SGNN to annotate dialog acts following clean up work to the original
# transcripts. Updates DAMSLTag (column I) and DialogAct (column J) values
# for the rows whose annotations changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 18;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 19;  Tag = "%";  Act = "Uninterpretable" },
    @{ Row = 27;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 36;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 40;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 44;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 46;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 47;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 48;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 51;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 70;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 72;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 81;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 89;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 90;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 104; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 115; Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 125; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 147; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 150; Tag = "ba"; Act = "Appreciation" },
    @{ Row = 155; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 163; Tag = "qy"; Act = "Yes-No-Question" },
    @{ Row = 168; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 175; Tag = "aa"; Act = "Agree/Accept" }
)

foreach ($u in $updates) {
    $ws.Range("I" + $u.Row).Value = $u.Tag
    $ws.Range("J" + $u.Row).Value = $u.Act
}
